# PELM Contact Members - add a new contact row (Agostinho Rocha) and
# move the active selection, mirroring the authored changes in the
# canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 currently holds empty, bordered placeholder cells (style "s=1").
# Reuse the existing "data row" formatting (border + centered Times/Calibri
# text, same as rows 3-4) by copying it down, rather than re-building the
# cell style from scratch - this keeps the resulting cellXfs identical to
# the ones already used by rows 3 and 4 instead of minting new ones.
$ws.Range("B3:E3").Copy() | Out-Null
$ws.Range("B5:E5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new contact's details.
$ws.Range("B5").Value = "Agostinho Rocha"
$ws.Range("C5").Value = "Portugal "
$ws.Range("D5").Value = "agostinho@fe.up.pt"
$ws.Range("E5").Value = 910563658

# Turn the e-mail address into a mailto hyperlink (adds the relationship +
# applies the built-in "Hyperlink" style/font on top of the border &
# centered alignment already pasted above).
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:agostinho@fe.up.pt") | Out-Null

# Move the active cell/selection to C6, matching the saved workbook view.
$ws.Range("C6").Select() | Out-Null
